# fix #683: replace the "geo" axis/column header by "country"
# (tutorial_IO notebook + doctests of read_csv/excel/hdf use "country" now)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("births")
$ws2 = $wb.Worksheets.Item("deaths")

$ws1.Range("A1").Value = "country"
$ws2.Range("A1").Value = "country"

# the "deaths" sheet had a stale selection (D17) left over from editing;
# reset the active cell back to A1 now that the sheet is clean
$ws2.Activate()
$ws2.Range("A1").Select()
